$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
# (matching the source inlineStr cells) need Text number format first,
# otherwise Excel auto-converts the assigned string into a real number.
$textCells = @("D5","D6","D7","D8","D10","D11","D12","D13","D15","D16","D18","D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D45","D46","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '70.385.31'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '3.550.52'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '617.68'
$ws.Range("E5").Value = '  +5.30%  '
$ws.Range("D6").Value = '187.32'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  +1.15%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").Value = '0.659'
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '53.82'
$ws.Range("E11").Value = '  -1.29%  '
$ws.Range("D12").Value = '0.0000310'
$ws.Range("E12").Value = '  -3.84%  '
$ws.Range("D13").Value = '9.69'
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("D14").Value = '4.110.60'
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").Value = '619.98'
$ws.Range("E15").Value = '  +8.16%  '
$ws.Range("D16").Value = '12.90'
$ws.Range("E16").Value = '  +3.83%  '
$ws.Range("D17").Value = '70.367.49'
$ws.Range("E17").Value = '  -0.90%  '
$ws.Range("D18").Value = '19.15'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '3.538.02'
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").Value = '17.75'
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").Value = '104.37'
$ws.Range("E23").Value = '  +9.72%  '
$ws.Range("D24").Value = '4.74'
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").Value = '5.12'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  +3.43%  '
$ws.Range("D27").Value = '11.01'
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").Value = '9.93'
$ws.Range("E28").Value = '  +8.72%  '
$ws.Range("D29").Value = '34.27'
$ws.Range("E29").Value = '  +5.84%  '
$ws.Range("D30").Value = '7.09'
$ws.Range("E30").Value = '  -1.86%  '
$ws.Range("D31").Value = '12.53'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("D33").Value = '64.27'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '3.70'
$ws.Range("E34").Value = '  +15.57%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").Value = '540.13'
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '3.17'
$ws.Range("E36").Value = '  -3.96%  '
$ws.Range("D37").Value = '0.401'
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Value = '37.38'
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("E40").Value = '  +3.76%  '
$ws.Range("D41").Value = '0.0₃0783'
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("D42").Value = '3.547.21'
$ws.Range("E42").Value = '  +1.46%  '
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("E44").Value = '  +5.80%  '
$ws.Range("D45").Value = '2.96'
$ws.Range("D46").Value = '0.144'
$ws.Range("E46").Value = '  +4.04%  '
$ws.Range("D47").Value = '3.39'
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("E48").Value = '  -3.86%  '
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("E50").Value = '  -3.26%  '
$ws.Range("D51").Value = '134.33'
$ws.Range("E51").Value = '  -0.43%  '
